$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column F (Actual Result) new values ----------------------------------
$ws.Range("F39").Value  = "Logo should be linked to the home page"
$ws.Range("F41").Value  = "Home page should be open after clicking on the logo"
$ws.Range("F43").Value  = "Home page should be open from every pages of the website after clicking on the logo"
$ws.Range("F45").Value  = "linked texts should be connected to the correct pages"
$ws.Range("F47").Value  = "Selected header links should be remain active"
$ws.Range("F49").Value  = "Shopping cart or checkout option sjould be open after clicking on cart icon"
$ws.Range("F51").Value  = "User profile should be shown after clicking on the profile icon"
$ws.Range("F53").Value  = "after clicking on the country flag icon user information should be shown by based on user credential"
$ws.Range("F55").Value  = "Correct information should be shown after clicking on the country logo"
$ws.Range("F57").Value  = "Another navigation part should be attached which will be showing country wise free shipping cost"
$ws.Range("F59").Value  = "Navigation bar should be functional"
$ws.Range("F61").Value  = "Header part should be attacvhed or sticky on the website"
$ws.Range("F64").Value  = "Footer should be perfectly aligned at the end of the website"
$ws.Range("F66").Value  = "There should be no extra white space from anyside at the end of the website "
$ws.Range("F68").Value  = "All the elements should be perfectly aligned in footer"
$ws.Range("F70").Value  = "There should be no grammatical mistakes in the footer"
$ws.Range("F72").Value  = "There should be no spelling mistakes"
$ws.Range("F74").Value  = "There should be equal space for the elements ( icons, linked texts, buttons)"
$ws.Range("F76").Value  = "Copyrights text should be added in the footer"
$ws.Range("F78").Value  = "There should be logo in the footer section"
$ws.Range("F80").Value  = "There should be a individual space for the logo in the footer"
$ws.Range("F82").Value  = "The  logo should be perfectly aligned"
$ws.Range("F84").Value  = "All the social icons should be perfectly aligned in there own space"
$ws.Range("F86").Value  = "There should be a pointer icon, by click on it , user will go to the top of the website from the bottom"
$ws.Range("F88").Value  = "There should be privacy policy, terms and use texts in the footer"
$ws.Range("F90").Value  = "Thereshould be a sign up option in the footer"
$ws.Range("F92").Value  = "Company number and eamil should be attached in the footer"
$ws.Range("F95").Value  = "All links of the footer should be working"
$ws.Range("F97").Value  = "Right page should be open after clicking on the linked texts"
$ws.Range("F99").Value  = "Home page should be open after clicking on the logo"
$ws.Range("F101").Value = "Social icons should be connected to the correct link and pages"
$ws.Range("F103").Value = "There should be no broken links"
$ws.Range("F105").Value = "There should be a subscription option"
$ws.Range("F107").Value = "Signup option shouldn't be taken invalid data"
$ws.Range("F109").Value = "Error message should be shpwn after giving invalid data"
$ws.Range("F111").Value = "Valid mail should be successfully submitted"
$ws.Range("F113").Value = "A pop-up or a message should be showen after successfully submitting the valid data"

# --- Style-only fixes on the blank spacer rows below two of the rows above -
# F84/F85 go from "left, no-wrap" (style 7) to "left, wrap" (style 2)
$ws.Range("F84:F85").WrapText = $true
# F103/F104 go from "center, wrap" (style 14) to "left, wrap" (style 2)
$ws.Range("F103:F104").HorizontalAlignment = -4131   # xlLeft

# --- Column D (Test case execution start date) -- added LAST so this shared
# string lands at the end of sharedStrings.xml (uniqueCount index 134), same
# order as the target diff.
$ws.Range("D3").Value = "15/10/2024"
$ws.Range("D3:D4").HorizontalAlignment = -4108   # xlCenter -> matches style 14

# --- View state: scroll position + selection -------------------------------
$excel.ActiveWindow.ScrollRow = 37
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D3:D4").Select()
